$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the sheet name (date rolled from 2021-11-05 to 2021-11-06)
$ws.Name = "Through 2021-11-06"

# Update the row label for November in column A (row 12)
$ws.Range("A12").Value = "November (through 11-06)"

# Update the November row (row 12) values
$ws.Range("B12").Value = 7
$ws.Range("C12").Value = 15
$ws.Range("D12").Value = 21
$ws.Range("E12").Value = 18
$ws.Range("F12").Value = 9
$ws.Range("G12").Value = 43
$ws.Range("H12").Value = 42

# Update the Total row (row 13) values
$ws.Range("B13").Value = 265
$ws.Range("C13").Value = 501
$ws.Range("D13").Value = 731
$ws.Range("E13").Value = 633
$ws.Range("F13").Value = 491
$ws.Range("G13").Value = 1100
$ws.Range("H13").Value = 1486
